$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 13 for the "Docentes responsaveis" value (was missing before) ---
# This shifts old rows 13-23 down to 14-24.
$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()

# Copy formatting (styles) for B13/C13 from an existing normal/red value-cell pair (row 2)
$ws.Range("B2:C2").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13").Value = '5840897 - Clodoaldo Saron'
$ws.Range("C13").Value = '5840897 - Clodoaldo Saron'

# --- Row 10 (Objetivos:) now gets the proper Portuguese objectives paragraph ---
$ws.Range("B10").Value = 'A reologia é a ciência que estuda o escoamento de materiais. O seu conhecimento é necessário para poder entender o processamento dos materiais poliméricos. A disciplina visa o ensino dos conceitos básicos de reologia de materiais (polímeros fundidos) para o estudante de engenharia de materiais, a disciplina visa também familiarizar o futuro engenheiro com os métodos experimentais para avaliação das propriedades reológicas de materiais poliméricos.Fornecer conhecimentos técnicos para o aluno escolher corretamente a técnica mais adequada de processamento de polímeros, bem como poder manipular e especificar corretamente as matérias primas, os equipamentos de processamento, os moldes e as máquinas adequadamente para determinadas conformações.'
$ws.Range("C10").Value = 'A reologia é a ciência que estuda o escoamento de materiais. O seu conhecimento é necessário para poder entender o processamento dos materiais poliméricos. A disciplina visa o ensino dos conceitos básicos de reologia de materiais (polímeros fundidos) para o estudante de engenharia de materiais, a disciplina visa também familiarizar o futuro engenheiro com os métodos experimentais para avaliação das propriedades reológicas de materiais poliméricos.Fornecer conhecimentos técnicos para o aluno escolher corretamente a técnica mais adequada de processamento de polímeros, bem como poder manipular e especificar corretamente as matérias primas, os equipamentos de processamento, os moldes e as máquinas adequadamente para determinadas conformações.'

# --- Row 14 (Programa resumido:, was row 13 before insert) gets the proper short-syllabus PT text ---
$ws.Range("B14").Value = 'Fundamentos de reologia. Processamento de polímeros: matérias-primas, máquinas e moldes.'
$ws.Range("C14").Value = 'Fundamentos de reologia. Processamento de polímeros: matérias-primas, máquinas e moldes.'

# --- Row 16 (Programa:, was row 15 before insert) gets the full PT syllabus text ---
$ws.Range("B16").Value = '1. Introdução a reologia. Tipos de fluxo. 2. Sólidos hookeanos e fluidos newtonianos. 3. Fluidos newtonianos e não newtonianos. 4. Viscoelasticidade. Viscosidade extensional. Diferenças de tensões normais. Variáveis que afetam a viscosidade de polímeros. 5. Importância da Reologia no processamento de polímeros. Fluxos utilizados para caracterizar materiais: fluxo de arraste, fluxos devido a diferença de pressão e escoamento em dutos. 6. Extrusão de polímeros: equipamentos, roscas, matrizes e aplicações. 7. Injeção de polímeros: equipamento, moldes, controle da operação, correção de problemas e aplicações. 8. Outras técnicas de processamento de termoplásticos: sopro, prensagem, termoformagem, calandragem, fiação, rotomoldagem. 9. Blendas e Compósitos Poliméricas: formas de obtenção, miscibilidade, compatibilidade e aplicações. 10. Técnicas de processamento de polímeros termorrígidos: moldagem manual, moldagem por pistola, pultrusão, enrolamento de filamento, prensagem, etc.'
$ws.Range("C16").Value = '1. Introdução a reologia. Tipos de fluxo. 2. Sólidos hookeanos e fluidos newtonianos. 3. Fluidos newtonianos e não newtonianos. 4. Viscoelasticidade. Viscosidade extensional. Diferenças de tensões normais. Variáveis que afetam a viscosidade de polímeros. 5. Importância da Reologia no processamento de polímeros. Fluxos utilizados para caracterizar materiais: fluxo de arraste, fluxos devido a diferença de pressão e escoamento em dutos. 6. Extrusão de polímeros: equipamentos, roscas, matrizes e aplicações. 7. Injeção de polímeros: equipamento, moldes, controle da operação, correção de problemas e aplicações. 8. Outras técnicas de processamento de termoplásticos: sopro, prensagem, termoformagem, calandragem, fiação, rotomoldagem. 9. Blendas e Compósitos Poliméricas: formas de obtenção, miscibilidade, compatibilidade e aplicações. 10. Técnicas de processamento de polímeros termorrígidos: moldagem manual, moldagem por pistola, pultrusão, enrolamento de filamento, prensagem, etc.'

# --- Row 19 (Metodo:, was row 18 before insert) gets the evaluation-method text ---
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'

# --- Row 20 (Criterio:, was row 19 before insert) gets the final-grade formula text ---
$ws.Range("B20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2'
$ws.Range("C20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2'

# --- Row 21 (Norma de recuperacao:, was row 20 before insert) gets the recovery-grade formula text ---
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'

# --- Row 22 (Bibliografia:, was row 21 before insert) gets the full bibliography text ---
$ws.Range("B22").Value = 'BRETAS, R. E. S.; D´ÁVILA, M. A. Reologia de Polímeros Fundidos, São Carlos, Eduscar, 2005.MANRICH, S. Processamento de termoplásticos – Rosca única, extrusão & matrizes, injeção & moldes,. McCRUM, N. G., BUCKLEY, C. P., BUCKNALl, C. B. Principles of Polymer Engineering, New York, Oxford University Press, 1997.Blass A., Processamento de Polímeros, editora da UFSC.CHAWLA, K. K. Composite Materials Science and Engineering, Spring-Verlag ed., Berlin, 1987.BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Livraria Medina, Coimbra, 1996.FONTANA, M. G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro, Ed. LTC, 2007 RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus, 1990SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2, 2000WOLYNEC, S. Técnicas Eletroquímicas em Corrosão, EDUSP, São Paulo, 2003'
$ws.Range("C22").Value = 'BRETAS, R. E. S.; D´ÁVILA, M. A. Reologia de Polímeros Fundidos, São Carlos, Eduscar, 2005.MANRICH, S. Processamento de termoplásticos – Rosca única, extrusão & matrizes, injeção & moldes,. McCRUM, N. G., BUCKLEY, C. P., BUCKNALl, C. B. Principles of Polymer Engineering, New York, Oxford University Press, 1997.Blass A., Processamento de Polímeros, editora da UFSC.CHAWLA, K. K. Composite Materials Science and Engineering, Spring-Verlag ed., Berlin, 1987.BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Livraria Medina, Coimbra, 1996.FONTANA, M. G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro, Ed. LTC, 2007 RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus, 1990SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2, 2000WOLYNEC, S. Técnicas Eletroquímicas em Corrosão, EDUSP, São Paulo, 2003'

